$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.612
$ws.Range("C3").Value = -12.285
$ws.Range("E3").Value = 16.548
$ws.Range("E6").Value = 16.71
$ws.Range("D8").Value = -8.401
$ws.Range("E10").Value = 16.661
$ws.Range("D11").Value = -7.391
$ws.Range("A12").Value = -21.629
$ws.Range("B14").Value = 6.228
$ws.Range("D14").Value = -7.528999999999999
$ws.Range("D15").Value = -8.175000000000001
$ws.Range("D17").Value = -8.446
$ws.Range("C20").Value = -12.149
$ws.Range("C25").Value = -12.431
$ws.Range("B26").Value = 6.142
$ws.Range("D26").Value = -8
$ws.Range("A27").Value = -21.589
$ws.Range("E27").Value = 16.561
$ws.Range("C30").Value = -12.717
$ws.Range("B31").Value = 6.371
$ws.Range("A32").Value = -21.439
$ws.Range("E33").Value = 17.558
$ws.Range("B35").Value = 7.657000000000001
$ws.Range("A36").Value = -20.761
$ws.Range("D36").Value = -8.375
$ws.Range("B37").Value = 7.641
$ws.Range("A38").Value = -20.223
$ws.Range("E39").Value = 16.643
$ws.Range("C44").Value = -12.646
$ws.Range("B45").Value = 5.718
$ws.Range("A46").Value = -21.515
$ws.Range("C47").Value = -12.316
$ws.Range("E47").Value = 16.335
$ws.Range("B52").Value = 5.4
$ws.Range("A54").Value = -21.646
$ws.Range("E54").Value = 16.554
$ws.Range("A55").Value = -21.899
$ws.Range("A56").Value = -21.818
$ws.Range("E56").Value = 16.653
$ws.Range("B57").Value = 6.090000000000001
$ws.Range("C58").Value = -12.98
$ws.Range("E58").Value = 16.675
$ws.Range("D64").Value = -7.81
$ws.Range("E66").Value = 17.114
$ws.Range("A67").Value = -21.588
$ws.Range("A69").Value = -21.604
$ws.Range("E69").Value = 17.13
$ws.Range("A72").Value = -21.436
$ws.Range("E72").Value = 16.603
$ws.Range("C78").Value = -12.536
$ws.Range("D79").Value = -7.813000000000001
$ws.Range("E80").Value = 16.44
$ws.Range("B81").Value = 6.303
$ws.Range("E82").Value = 16.936
$ws.Range("A83").Value = -21.1
$ws.Range("B83").Value = 6.77
$ws.Range("E83").Value = 16.694
$ws.Range("C84").Value = -13.001
$ws.Range("A86").Value = -22.216
$ws.Range("C89").Value = -12.116
$ws.Range("D89").Value = -7.423999999999999
$ws.Range("A91").Value = -21.632
$ws.Range("C91").Value = -11.069
$ws.Range("C92").Value = -11.509
$ws.Range("A93").Value = -21.49
$ws.Range("C96").Value = -13.02
$ws.Range("A99").Value = -20.682
$ws.Range("B100").Value = 5.789
$ws.Range("B102").Value = 7.039
$ws.Range("C102").Value = -12.779
